# Budget breakdown template: insert a new "Release" column before the
# existing "Policy Amount" column (old column G), shifting the policy
# amount / budget id / budget plan id columns one to the right, and
# tighten up a few row heights / column widths to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before G. Everything from G rightwards
#    (Policy Amount, the hidden budget_id/budget_plan_id helper columns,
#    and the trailing default-width columns) shifts one column to the
#    right; formulas such as B8's SUM(G12:G5000) get rewritten to
#    reference the new column automatically.
$ws.Columns("G:G").Insert()

# 2) Give the header row its new "Release" / "Latest Policy Amount"
#    labels: F11 used to read "Latest Policy Amount" (now "Release"),
#    and the freshly inserted, still-blank G11 becomes the old F11 text.
$ws.Range("G11").Value = $ws.Range("F11").Text
$ws.Range("F11").Value = "Release"

# 3) Column width / style tidy-up: the new column G keeps the wide
#    "Policy/Latest Policy Amount" formatting that used to live in F,
#    while F itself becomes narrower now that it's just "Release".
$ws.Columns("G:G").NumberFormat = "#,##0.00_);(#,##0.00)"
$ws.Columns("G:G").ColumnWidth = 24.22
$ws.Columns("F:F").ColumnWidth = 17.44

# 4) A handful of rows shrink slightly (12.85 -> 12.8, 13.5 -> 13.25).
$ws.Rows("5:10").RowHeight = 12.8
$ws.Rows("11:11").RowHeight = 13.25
